$d = $word.ActiveDocument

# Fix missing "ich" in the User Story-014 description:
# "Als Betreuer möchte mich im System anmelden können, um die Sicherheit zu erhöhen."
# -> "Als Betreuer möchte ich mich im System anmelden können, um die Sicherheit zu erhöhen."
$d.Content.Find.Execute(
    "Als Betreuer möchte mich im System anmelden können",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Als Betreuer möchte ich mich im System anmelden können",
    2)
